$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mineralization parameter (DOC_miner_const) value in B21
$ws.Range("B21").Value = 0.003

# Update the active cell selection to match the saved state
$ws.Range("B22").Select()
